# Update countries & provincias Spain
#
# This script applies the refreshed COVID-19 dataset values to the
# "Pais" worksheet, including:
#  - Updated "Datos actualizados..." timestamp (A1)
#  - Updated numeric counters for several countries (rank unaffected
#    except for the three-way/two-way reorders below)
#  - Two label swaps caused by the source data being re-sorted:
#      * Georgia now ranks above Republica de Chipre (rows 146/147)
#      * Montserrat now ranks above Islas Malvinas (rows 213/214)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp update (A1) ---------------------------------------------
$ws.Range("A1").Value() = "Datos actualizados a 16 de Agosto de 2020 a las 09:44"

# --- Row 7: Rusia --------------------------------------------------------
$ws.Range("B7").Value() = 922853
$ws.Range("C7").Value() = 4969
$ws.Range("D7").Value() = 732968
$ws.Range("E7").Value() = 174200
$ws.Range("F7").Value() = 0
$ws.Range("G7").Value() = 68
$ws.Range("H7").Value() = 15685

# --- Row 56: Armenia -------------------------------------------------------
$ws.Range("B56").Value() = 41663
$ws.Range("C56").Value() = 168
$ws.Range("D56").Value() = 34584
$ws.Range("E56").Value() = 6261
$ws.Range("F56").Value() = 0
$ws.Range("G56").Value() = 1
$ws.Range("H56").Value() = 818

# --- Row 84: Sudan -----------------------------------------------------
$ws.Range("B84").Value() = 12314
$ws.Range("C84").Value() = 103
$ws.Range("D84").Value() = 6350
$ws.Range("E84").Value() = 5166
$ws.Range("F84").Value() = 0
$ws.Range("G84").Value() = 2
$ws.Range("H84").Value() = 798

# --- Row 108: Hungria ----------------------------------------------------
$ws.Range("B108").Value() = 4916
$ws.Range("C108").Value() = 39
$ws.Range("D108").Value() = 3623
$ws.Range("E108").Value() = 685
$ws.Range("F108").Value() = 0
$ws.Range("G108").Value() = 1
$ws.Range("H108").Value() = 608

# --- Rows 146/147: Georgia now outranks Republica de Chipre --------------
# Row 146 becomes Georgia (updated figures)
$ws.Range("A146").Value() = "Georgia"
$ws.Range("B146").Value() = 1336
$ws.Range("C146").Value() = 15
$ws.Range("D146").Value() = 1088
$ws.Range("E146").Value() = 231
$ws.Range("F146").Value() = 0
$ws.Range("G146").Value() = 0
$ws.Range("H146").Value() = 17

# Row 147 becomes Republica de Chipre (unchanged figures, just moved down)
$ws.Range("A147").Value() = "Republica de Chipre"
$ws.Range("B147").Value() = 1332
$ws.Range("C147").Value() = 0
$ws.Range("D147").Value() = 870
$ws.Range("E147").Value() = 442
$ws.Range("F147").Value() = 0
$ws.Range("G147").Value() = 0
$ws.Range("H147").Value() = 20

# --- Row 148: Letonia ------------------------------------------------------
$ws.Range("B148").Value() = 1322
$ws.Range("C148").Value() = 7
$ws.Range("D148").Value() = 1078
$ws.Range("E148").Value() = 212
$ws.Range("F148").Value() = 0
$ws.Range("G148").Value() = 0
$ws.Range("H148").Value() = 32

# --- Row 155: Jamaica ------------------------------------------------------
$ws.Range("B155").Value() = 1106
$ws.Range("C155").Value() = 24
$ws.Range("D155").Value() = 761
$ws.Range("E155").Value() = 331
$ws.Range("F155").Value() = 0
$ws.Range("G155").Value() = 0
$ws.Range("H155").Value() = 14

# --- Row 199: Curazao -----------------------------------------------------
$ws.Range("B199").Value() = 34
$ws.Range("C199").Value() = 1
$ws.Range("D199").Value() = 31
$ws.Range("E199").Value() = 2
$ws.Range("F199").Value() = 0
$ws.Range("G199").Value() = 0
$ws.Range("H199").Value() = 1

# --- Rows 213/214: Montserrat now outranks Islas Malvinas -----------------
# Row 213 becomes Montserrat (unchanged figures, just moved up)
$ws.Range("A213").Value() = "Montserrat"
$ws.Range("B213").Value() = 13
$ws.Range("C213").Value() = 0
$ws.Range("D213").Value() = 12
$ws.Range("E213").Value() = 0
$ws.Range("F213").Value() = 0
$ws.Range("G213").Value() = 0
$ws.Range("H213").Value() = 1

# Row 214 becomes Islas Malvinas (unchanged figures, just moved down)
$ws.Range("A214").Value() = "Islas Malvinas"
$ws.Range("B214").Value() = 13
$ws.Range("C214").Value() = 0
$ws.Range("D214").Value() = 13
$ws.Range("E214").Value() = 0
$ws.Range("F214").Value() = 0
$ws.Range("G214").Value() = 0
$ws.Range("H214").Value() = 0
